$d = $word.ActiveDocument

# Locate the unique "<add>r</add>, " span (literal text, as the tags are
# stored as plain visible text runs in this markup-annotated document).
$rng = $d.Content
$rng.Find.Execute("<add>r</add>, ", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "", 0)

if (-not $rng.Find.Found) {
    throw "Target text '<add>r</add>, ' not found"
}

$tagStart = $rng.Start

# Remove the trailing </add> tag run first (so earlier offsets stay valid),
# then the leading <add> tag run. This leaves the existing "r" and ", "
# runs (which already carry the plain color=000000 formatting) adjacent
# to each other, so Word merges them into a single "r, " run.
$addClose = $d.Range($tagStart + 6, $tagStart + 12)
$addClose.Delete()

$addOpen = $d.Range($tagStart, $tagStart + 5)
$addOpen.Delete()


